# Update Planung/TODOS.xlsx (Tabelle1) per commit:
# "Diverse Planungs-Dokumente aktualisiert, Vorlage Aspektmigration"

$xlPasteFormats = -4122  # XlPasteType.xlPasteFormats

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Fix typo "Auwände" -> "Aufwände" in the task description (row 17)
$ws.Range("A17").Value = "Aufwände anpassen, Präsentationstechnikkurs"

# Mark several items as done ("Erledigt") in the Status column (C)
# Rows 8, 9, 11 already have a Status cell, but it was styled for the old
# "To be reviewed"/blank state (red font). Re-format to the plain "Erledigt"
# look (copy from C2) before writing the new value.
$ws.Range("C2").Copy()
$ws.Range("C8").PasteSpecial($xlPasteFormats)
$ws.Range("C8").Value = "Erledigt"

$ws.Range("C2").Copy()
$ws.Range("C9").PasteSpecial($xlPasteFormats)
$ws.Range("C9").Value = "Erledigt"

$ws.Range("C2").Copy()
$ws.Range("C11").PasteSpecial($xlPasteFormats)
$ws.Range("C11").Value = "Erledigt"

# Rows 17, 18, 19 get a brand-new Status cell; copy formatting from an existing
# "Status = Erledigt" cell (C2) so borders/font match the rest of the column.
$ws.Range("C2").Copy()
$ws.Range("C17").PasteSpecial($xlPasteFormats)
$ws.Range("C17").Value = "Erledigt"

$ws.Range("C2").Copy()
$ws.Range("C18").PasteSpecial($xlPasteFormats)
$ws.Range("C18").Value = "Erledigt"

$ws.Range("C2").Copy()
$ws.Range("C19").PasteSpecial($xlPasteFormats)
$ws.Range("C19").Value = "Erledigt"

# Add new "Termin" (deadline) dates in column D for a few rows.
# Copy formatting from existing date cells so number format/border/font match.

# D6 / D27 -> 01.06.2013 (normal style, like D29/D30/D31/D32)
$ws.Range("D29").Copy()
$ws.Range("D6").PasteSpecial($xlPasteFormats)
$ws.Range("D6").Value = 41426

$ws.Range("D29").Copy()
$ws.Range("D27").PasteSpecial($xlPasteFormats)
$ws.Range("D27").Value = 41426

# D23 / D25 -> 21.05.2013 (bold red "urgent" style, like D28)
$ws.Range("D28").Copy()
$ws.Range("D23").PasteSpecial($xlPasteFormats)
$ws.Range("D23").Value = 41415

$ws.Range("D28").Copy()
$ws.Range("D25").PasteSpecial($xlPasteFormats)
$ws.Range("D25").Value = 41415

# Move the active selection to D33 (matches the saved view state in the file)
[void]$ws.Range("D33").Select()

Write-Host "edits applied"
